# Applies the "quests export" change:
#   - Adds 5 new NPC rows to the "NPCs" sheet (rows 11-15)
#   - Adds 5 matching rows to the "Npcs Commands" sheet (rows 11-15)
#   - Widens column D on the "NPCs" sheet to fit the new, longer
#     game_map_id value ("Shadow Plane")

$wb = $excel.ActiveWorkbook

$npcs = $wb.Worksheets.Item("NPCs")
$cmds = $wb.Worksheets.Item("Npcs Commands")

# ---------------------------------------------------------------------
# NPCs sheet - new rows 11-15
# Columns: A name | B real_name | C type | D game_map_id |
#          E moves_around_map | F must_be_at_same_location |
#          G text_command_to_message | H x_position | I y_position
# ---------------------------------------------------------------------

$npcs.Cells.Item(11, 1).Value = "ShadeLord"
$npcs.Cells.Item(11, 2).Value = "Shade Lord"
$npcs.Cells.Item(11, 3).Value = 2
$npcs.Cells.Item(11, 4).Value = "Shadow Plane"
$npcs.Cells.Item(11, 6).Value = 1
$npcs.Cells.Item(11, 7).Value = "/m ShadeLord:"
$npcs.Cells.Item(11, 8).Value = 288
$npcs.Cells.Item(11, 9).Value = 480

$npcs.Cells.Item(12, 1).Value = "DrunkenAdventurer"
$npcs.Cells.Item(12, 2).Value = "DrunkenAdventurer"
$npcs.Cells.Item(12, 3).Value = 2
$npcs.Cells.Item(12, 4).Value = "Dungeons"
$npcs.Cells.Item(12, 6).Value = 1
$npcs.Cells.Item(12, 7).Value = "/m DrunkenAdventurer:"
$npcs.Cells.Item(12, 8).Value = 384
$npcs.Cells.Item(12, 9).Value = 368

$npcs.Cells.Item(13, 1).Value = "ChildofShade"
$npcs.Cells.Item(13, 2).Value = "Child of Shade"
$npcs.Cells.Item(13, 3).Value = 1
$npcs.Cells.Item(13, 4).Value = "Shadow Plane"
$npcs.Cells.Item(13, 6).Value = 1
$npcs.Cells.Item(13, 7).Value = "/m ChildofShade:"
$npcs.Cells.Item(13, 8).Value = 320
$npcs.Cells.Item(13, 9).Value = 288

$npcs.Cells.Item(14, 1).Value = "HellsGateKeeper"
$npcs.Cells.Item(14, 2).Value = "HellsGateKeeper"
$npcs.Cells.Item(14, 3).Value = 2
$npcs.Cells.Item(14, 4).Value = "Hell"
$npcs.Cells.Item(14, 6).Value = 1
$npcs.Cells.Item(14, 7).Value = "/m HellsGateKeeper:"
$npcs.Cells.Item(14, 8).Value = 16
$npcs.Cells.Item(14, 9).Value = 336

$npcs.Cells.Item(15, 1).Value = "QueenofHearts"
$npcs.Cells.Item(15, 2).Value = "Queen of Hearts"
$npcs.Cells.Item(15, 3).Value = 3
$npcs.Cells.Item(15, 4).Value = "Hell"
$npcs.Cells.Item(15, 6).Value = 1
$npcs.Cells.Item(15, 7).Value = "/m QueenofHearts:"
$npcs.Cells.Item(15, 8).Value = 432
$npcs.Cells.Item(15, 9).Value = 80

# Column D needs to be a bit wider now that "Shadow Plane" lives there.
$npcs.Range("D1").ColumnWidth = 14.1666666666667

# ---------------------------------------------------------------------
# Npcs Commands sheet - new rows 11-15
# Columns: A npc_id | B command | C command_type
# ---------------------------------------------------------------------

$cmds.Cells.Item(11, 1).Value = "Shade Lord"
$cmds.Cells.Item(11, 2).Value = "Shades"
$cmds.Cells.Item(11, 3).Value = 0

$cmds.Cells.Item(12, 1).Value = "DrunkenAdventurer"
$cmds.Cells.Item(12, 2).Value = "Story"
$cmds.Cells.Item(12, 3).Value = 0

$cmds.Cells.Item(13, 1).Value = "Child of Shade"
$cmds.Cells.Item(13, 2).Value = "Shadows"
$cmds.Cells.Item(13, 3).Value = 2

$cmds.Cells.Item(14, 1).Value = "HellsGateKeeper"
$cmds.Cells.Item(14, 2).Value = "Gates"
$cmds.Cells.Item(14, 3).Value = 0

$cmds.Cells.Item(15, 1).Value = "Queen of Hearts"
$cmds.Cells.Item(15, 2).Value = "Enchantments"
$cmds.Cells.Item(15, 3).Value = 3
